$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (C) column holds the last-changed date for every entry.
# This automatic update bumps that date forward by one day (45243 -> 45244)
# for every data row in the sheet.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
